$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "30.005.50"
$cell.ClearFormats()
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  -0.26%  "
$cell.ClearFormats()
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.870.72"
$cell.ClearFormats()
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  -2.49%  "
$cell.ClearFormats()
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.ClearFormats()
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  -0.04%  "
$cell.ClearFormats()
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "319.19"
$cell.ClearFormats()
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  -3.35%  "
$cell.ClearFormats()
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  +0.00%  "
$cell.ClearFormats()
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.5046"
$cell.ClearFormats()
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  -2.94%  "
$cell.ClearFormats()
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3965"
$cell.ClearFormats()
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  -2.81%  "
$cell.ClearFormats()
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.08201"
$cell.ClearFormats()
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  -3.50%  "
$cell.ClearFormats()
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "42.11"
$cell.ClearFormats()
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  -2.47%  "
$cell.ClearFormats()
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  -3.10%  "
$cell.ClearFormats()
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "23.42"
$cell.ClearFormats()
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  +4.42%  "
$cell.ClearFormats()
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.876.81"
$cell.ClearFormats()
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  -3.64%  "
$cell.ClearFormats()
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.289"
$cell.ClearFormats()
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  -1.90%  "
$cell.ClearFormats()
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "7.189"
$cell.ClearFormats()
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  -2.87%  "
$cell.ClearFormats()
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.ClearFormats()
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  +0.04%  "
$cell.ClearFormats()
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "91.87"
$cell.ClearFormats()
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  -3.86%  "
$cell.ClearFormats()
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -2.31%  "
$cell.ClearFormats()
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06423"
$cell.ClearFormats()
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  -0.88%  "
$cell.ClearFormats()
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  +0.05%  "
$cell.ClearFormats()
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "29.994.04"
$cell.ClearFormats()
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  -0.36%  "
$cell.ClearFormats()
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.851"
$cell.ClearFormats()
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  -2.60%  "
$cell.ClearFormats()
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "11.14"
$cell.ClearFormats()
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  -1.51%  "
$cell.ClearFormats()
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.169"
$cell.ClearFormats()
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -2.62%  "
$cell.ClearFormats()
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.088.31"
$cell.ClearFormats()
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  -2.71%  "
$cell.ClearFormats()
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "21.22"
$cell.ClearFormats()
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  +0.56%  "
$cell.ClearFormats()
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "160.94"
$cell.ClearFormats()
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  +0.45%  "
$cell.ClearFormats()
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.216"
$cell.ClearFormats()
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  -9.46%  "
$cell.ClearFormats()
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "127.22"
$cell.ClearFormats()
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -1.45%  "
$cell.ClearFormats()
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.070"
$cell.ClearFormats()
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  -0.43%  "
$cell.ClearFormats()
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  -1.99%  "
$cell.ClearFormats()
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  -2.31%  "
$cell.ClearFormats()
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "3.648"
$cell.ClearFormats()
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  +0.35%  "
$cell.ClearFormats()
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.02439"
$cell.ClearFormats()
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  -2.15%  "
$cell.ClearFormats()
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "5.213"
$cell.ClearFormats()
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  +0.53%  "
$cell.ClearFormats()
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.06358"
$cell.ClearFormats()
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -3.91%  "
$cell.ClearFormats()
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.2141"
$cell.ClearFormats()
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -2.94%  "
$cell.ClearFormats()
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.173"
$cell.ClearFormats()
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  -4.60%  "
$cell.ClearFormats()
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "8.485"
$cell.ClearFormats()
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  -4.76%  "
$cell.ClearFormats()
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  -2.41%  "
$cell.ClearFormats()
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.6303"
$cell.ClearFormats()
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  -3.15%  "
$cell.ClearFormats()
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "11.28"
$cell.ClearFormats()
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  -2.82%  "
$cell.ClearFormats()
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.ClearFormats()
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  -0.12%  "
$cell.ClearFormats()
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "12.99"
$cell.ClearFormats()
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  -2.17%  "
$cell.ClearFormats()
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.5914"
$cell.ClearFormats()
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  -3.98%  "
$cell.ClearFormats()
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.088"
$cell.ClearFormats()
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  +0.34%  "
$cell.ClearFormats()
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "3.623"
$cell.ClearFormats()
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -3.95%  "
$cell.ClearFormats()
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "122.52"
$cell.ClearFormats()
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  -1.50%  "
$cell.ClearFormats()
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.206"
$cell.ClearFormats()
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  -3.16%  "
$cell.ClearFormats()
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "77.48"
$cell.ClearFormats()
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  -2.95%  "
$cell.ClearFormats()
